# Fix: import was missing 2 fields (Highlight & Chuong trinh dac biet) and
# the product-detail view was missing url3/url4 columns.
#
# This inserts 4 new header columns (O1:R1) into row 1 of the import
# template, pushing the former O1:T1 headers (Thuong hieu, Cong nghe,
# So kenh, Phan giai, Luu tru, Tinh nang) out to S1:X1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: Highlight, Chuong trinh dac biet, and two more promo columns
# (re-using the existing "Khuyen mai" header text, matching the source file).
$ws.Range("O1").Value = "Highlight"
$ws.Range("P1").Value = "Chương trình đặc biệt"
$ws.Range("Q1").Value = "Khuyến mại"
$ws.Range("R1").Value = "Khuyến mại"

# Re-write the headers that shifted right into their new home, S1:X1
# (literal values taken from the original O1:T1, since COM Range.Value
# read-back isn't reliable in this host -- write the known text instead).
$ws.Range("S1").Value = "Thương hiệu"
$ws.Range("T1").Value = "Công nghệ"
$ws.Range("U1").Value = "Số kênh"
$ws.Range("V1").Value = "Phân giải"
$ws.Range("W1").Value = "Lưu trữ"
$ws.Range("X1").Value = "Tính năng"

# U1:X1 are brand-new cells beyond the sheet's old T-column extent, so they
# would otherwise inherit the plain default column formatting instead of
# the bold Times New Roman header style used by the rest of row 1. Force
# the header font explicitly so they match the "s=1" header style.
$headerRange = $ws.Range("U1:X1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 13
$headerRange.Font.Name = "Times New Roman"
